$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B17: currently 60, becomes formula 50+40 = 90
$ws.Range("B17").Formula = "=50 + 40"

# Update C17 shared string "init backend" -> "init backend, modularisointi ja mongoose alustus"
$ws.Range("C17").Value = "init backend, modularisointi ja mongoose alustus"

# Update selection from D17 to C18
$ws.Range("C18").Select()
